# Update latest output (run 13)
$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update cost/unit cost for the single schedule row ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = 1686.649185
$schedule.Range("F2").Value = 27.88771800595238

# --- Sheet "Detailed": update price values (and historical/forecast labels) ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B21").Value = 48.97284
$detailed.Range("C21").Value = "historical"

$detailed.Range("B22").Value = 36.06
$detailed.Range("C22").Value = "historical"

$detailed.Range("B23").Value = 36.06

$detailed.Range("B24").Value = 34.01

$detailed.Range("B25").Value = 36.06

$detailed.Range("B26").Value = 36.06

$detailed.Range("B28").Value = 11.92153

$detailed.Range("B29").Value = 0.51

$detailed.Range("B30").Value = 49.7961

$detailed.Range("B32").Value = 36.01246

$detailed.Range("B33").Value = 36.01246

$detailed.Range("B34").Value = 40.69742

$detailed.Range("B35").Value = 41.04596

$detailed.Range("B36").Value = 54.74532

$detailed.Range("B37").Value = 25.11183

$detailed.Range("B38").Value = 55.11546

$detailed.Range("B39").Value = 72.95139

$detailed.Range("B41").Value = 158.99

$detailed.Range("B42").Value = 158.99

$detailed.Range("B44").Value = 105.79

$detailed.Range("B45").Value = 85.95

$detailed.Range("B46").Value = 71.40000000000001

$detailed.Range("B47").Value = 64.99985

$detailed.Range("B48").Value = 70.27191000000001

$detailed.Range("B49").Value = 60.18313
